# Insert a new weekly record row right before the existing row 454
# (Feria Lagunitas de Puerto Montt - Zanahoria), shifting all subsequent
# rows (old 454..515) down by one (to 455..516).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 454; Excel shifts rows 454..515 down
# to 455..516 automatically.
$ws.Rows.Item(454).Insert()

# Populate the newly inserted row 454 with the new weekly record.
$ws.Cells.Item(454, 1).Value = 4
$ws.Cells.Item(454, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(454, 3).Value = "Los Lagos"
$ws.Cells.Item(454, 4).Value = 44984
$ws.Cells.Item(454, 5).Value = 10
$ws.Cells.Item(454, 6).Value = 100114013
$ws.Cells.Item(454, 7).Value = "Zanahoria"
$ws.Cells.Item(454, 8).Value = "Sin especificar"
$ws.Cells.Item(454, 9).Value = "Primera"
$ws.Cells.Item(454, 10).Value = 150
$ws.Cells.Item(454, 11).Value = 10000
$ws.Cells.Item(454, 12).Value = 10000
$ws.Cells.Item(454, 13).Value = 10000
$ws.Cells.Item(454, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(454, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(454, 16).Value = 500
$ws.Cells.Item(454, 17).Value = 20
$ws.Cells.Item(454, 18).Value = "Hortaliza"
